$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("study")
Write-Output ("Col17 Width=" + $ws.Columns.Item(17).Width)
